# participation.xlsx — "Add files via upload" edit
#
# Net effect (per the target diff): the single template row on Sheet1
# (Name | Event | Email-hyperlink | Subject) is reworked into
# (Name | Event | Subject | Email-hyperlink) — i.e. columns C and D are
# swapped — and that row is then replicated down through row 10 so every
# participant row (1-10) carries the same Name/Event/Subject/Email-hyperlink
# quadruple, each with its own hyperlink relationship.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Swap columns C and D -------------------------------------------
# Before: C = hyperlinked e-mail (style "Hyperlink"), D = plain "Computer"
# After:  C = plain "Computer",  D = hyperlinked e-mail
$ws.Columns.Item(4).Cut()
$ws.Columns.Item(3).Insert()

# Stash the canonical hyperlink-cell format (now sitting on D1) in a
# scratch cell so we can re-apply it later after Hyperlinks.Add touches it.
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Replicate the header/template row down to row 10 ---------------
$ws.Range("A1:D1").Copy()
for ($r = 2; $r -le 10; $r++) {
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4104)   # xlPasteAll
}
$excel.CutCopyMode = 0

# xlPasteAll approximates rather than exactly clones D1's xf onto D2:D10
# (it drops the numFmtId=14 quirk that rides along with the hyperlink
# font), so normalize the whole column back to the stashed canonical
# style before touching hyperlinks - that way every cell starts from the
# same format and Hyperlinks.Add below only has to mint one throw-away
# style variant instead of one per distinct starting format.
$ws.Range("F1").Copy()
$ws.Range("D1:D10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3. Fix up the hyperlink -------------------------------------------
# The column swap leaves the original Hyperlink object's anchor on the old
# (now plain) C1 cell instead of moving with the value to D1, so drop it
# and recreate one hyperlink per row, each pointing at the e-mail address
# already sitting in that row's D cell.
$ws.Range("C1").Hyperlinks.Delete()

for ($r = 1; $r -le 10; $r++) {
    $cell = $ws.Range("D" + $r)
    $ws.Hyperlinks.Add($cell, "mailto:halgodeshivraj03@gmail.com")
}

# Re-apply the canonical hyperlink style (clobbered by Hyperlinks.Add) to
# the whole D1:D10 column in one shot, then drop the scratch cell.
$ws.Range("F1").Copy()
$ws.Range("D1:D10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F1").Clear()
$excel.CutCopyMode = 0

# --- 4. Drop the now-unused trailing placeholder cell -------------------
$ws.Range("E1").Clear()

# --- 5. Selection, to match the saved view ------------------------------
$ws.Range("D12").Select()
